$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B4 = 1 (task "پنل کاربر" is now done)
$ws.Range("B4").Value = 1

# Add new task rows
$ws.Range("A6").Value = "صفحه سبد خرید"
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = "صغحه ویرایش اطلاعات کاربر"
$ws.Range("B7").Value = 1

# Set B9 = 1 (task "صفحه دسته بندی ها" is now done)
$ws.Range("B9").Value = 1

# Update the active selection to B3
$ws.Range("B3").Select() | Out-Null
